# Dividend Calculation.xlsx — update June 2017 Taxable Account dividend
# (all dependent SUM/shared formulas on both sheets recalculate automatically).

$wb = $excel.ActiveWorkbook

# --- "Yearly" sheet: bump the Taxable Account value for June 2017 (row 8, col L) ---
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsYearly.Range("L8").Value = 152.4

# Match the author's cursor/selection left behind in the saved file.
$wsYearly.Range("L9").Select()

# --- "All Time" sheet: depends on Yearly!L15 via formulas, recalculates on its own ---
$wsAllTime = $wb.Worksheets.Item("All Time")
$wsAllTime.Activate()
$wsAllTime.Range("F67").Select()

$wb.Save()
